# Apply changes described by the diff:
# - H2:I9 cells hold numbers that were stored as text (inlineStr); convert
#   them to real numeric values while keeping the same displayed value.
# - C10 text changes from "Run()" to "1.0"
# - H10 is cleared (becomes blank) while I10 becomes a real number (30)
# - I11 becomes a real number (60)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-9: H column and I column values, set as numbers.
$values = @{
    2 = @{ H = 2; I = 45 }
    3 = @{ H = 2; I = 30 }
    4 = @{ H = 2; I = 60 }
    5 = @{ H = 2; I = 60 }
    6 = @{ H = 2; I = 45 }
    7 = @{ H = 2; I = 45 }
    8 = @{ H = 1; I = 45 }
    9 = @{ H = 1; I = 45 }
}

foreach ($row in $values.Keys) {
    $ws.Range("H$row").Value = $values[$row].H
    $ws.Range("I$row").Value = $values[$row].I
}

# Row 10: C10 text updated (forced to stay text, not become the number 1),
# H10 cleared, I10 becomes numeric
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "1.0"
$ws.Range("C10").Style = "Normal"
$ws.Range("H10").ClearContents()
$ws.Range("I10").Value = 30

# Row 11: I11 becomes numeric
$ws.Range("I11").Value = 60
